$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins, Losses, Ties) in AD1:AF1,
# matching the formatting of the existing header row (bold font,
# thin border, centered horizontal / top vertical alignment).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Fill the season record values (Wins=95, Losses=67, Ties=0) for every
# player row (rows 2 through 47).
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 95
    $ws.Cells.Item($r, 31).Value = 67
    $ws.Cells.Item($r, 32).Value = 0
}
